$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at row 259 (header in row 1). Find the last
# populated row in column A so new rows are appended right after it.
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

# New monthly sales figures for month 7 (July) / year 2025, one row per store.
$newData = @(
    @(7, 1, 6805.15, 2025, "Bibi Cell Mundi"),
    @(7, 2, 2251,    2025, "Bibi Cell Manauara"),
    @(7, 3, 3638,    2025, "Bibi Cell Vieiralves"),
    @(7, 4, 4535.01, 2025, "Bibi Cell Ponta Negra")
)

foreach ($row in $newData) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $row[0]
    $ws.Cells.Item($lastRow, 2).Value = $row[1]
    $ws.Cells.Item($lastRow, 3).Value = $row[2]
    $ws.Cells.Item($lastRow, 4).Value = $row[3]
    $ws.Cells.Item($lastRow, 5).Value = $row[4]
}
